$wb = $excel.ActiveWorkbook

# "Ready for handoff" -> "In Translation" everywhere it appears (Overview!E:F,
# zh-cn!C, de-de!C all point at the same shared string).
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation") | Out-Null
}

# The status text got shorter, so the (autofit) Status-ish columns narrow too:
# Overview columns E ("zh-cn") and F ("de-de"), and column C ("Status") on the
# zh-cn / de-de detail sheets.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
